# Updates cryptos list values (prices in column D, volume % in column E)
# to match the refreshed data, plus a swap of the BabyDogeCoin/RenderToken rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.007.75'
$ws.Range("E2").Value = '  -4.06%  '
$ws.Range("D3").Value = '1.742.62'
$ws.Range("E3").Value = '  -4.53%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '''226.45'
$ws.Range("E5").Value = '  -3.51%  '
$ws.Range("D6").Value = '''0.5790'
$ws.Range("E6").Value = '  -3.58%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '''0.2724'
$ws.Range("E8").Value = '  -1.22%  '
$ws.Range("E9").Value = '  -1.17%  '
$ws.Range("D10").Value = '''0.06604'
$ws.Range("E10").Value = '  -4.77%  '
$ws.Range("D11").Value = '''0.07556'
$ws.Range("E11").Value = '  -0.56%  '
$ws.Range("D12").Value = '1.746.96'
$ws.Range("E12").Value = '  -4.58%  '
$ws.Range("D13").Value = '''4.712'
$ws.Range("E13").Value = '  -0.58%  '
$ws.Range("D14").Value = '''0.6033'
$ws.Range("E14").Value = '  -3.73%  '
$ws.Range("D15").Value = '1.979.49'
$ws.Range("E15").Value = '  -4.64%  '
$ws.Range("D16").Value = '''74.52'
$ws.Range("E16").Value = '  -3.83%  '
$ws.Range("D17").Value = '''0.000008738'
$ws.Range("E17").Value = '  -11.03%  '
$ws.Range("D18").Value = '28.012.02'
$ws.Range("E18").Value = '  -2.23%  '
$ws.Range("E19").Value = '  -4.02%  '
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").Value = '''205.47'
$ws.Range("E21").Value = '  -4.97%  '
$ws.Range("D22").Value = '''11.29'
$ws.Range("E22").Value = '  -2.23%  '
$ws.Range("D23").Value = '''6.628'
$ws.Range("E23").Value = '  -3.60%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").Value = '''150.04'
$ws.Range("E25").Value = '  -3.65%  '
$ws.Range("D26").Value = '''8.109'
$ws.Range("E26").Value = '  +2.26%  '
$ws.Range("D27").Value = '''0.1233'
$ws.Range("E27").Value = '  -4.26%  '
$ws.Range("D28").Value = '''16.16'
$ws.Range("E28").Value = '  -1.73%  '
$ws.Range("D29").Value = '''1.387'
$ws.Range("E29").Value = '  -1.85%  '
$ws.Range("D30").Value = '''0.06140'
$ws.Range("E30").Value = '  -4.84%  '
$ws.Range("D31").Value = '''1.393'
$ws.Range("E31").Value = '  -3.27%  '
$ws.Range("D32").Value = '''3.744'
$ws.Range("E32").Value = '  -2.17%  '
$ws.Range("D33").Value = '''3.721'
$ws.Range("E33").Value = '  -1.39%  '
$ws.Range("D34").Value = '''1.667'
$ws.Range("E34").Value = '  -3.33%  '
$ws.Range("D35").Value = '''1.037'
$ws.Range("E35").Value = '  -5.04%  '
$ws.Range("D36").Value = '''0.6395'
$ws.Range("E36").Value = '  -1.19%  '
$ws.Range("D37").Value = '''2.411'
$ws.Range("E37").Value = '  -4.98%  '
$ws.Range("D38").Value = '''2.713'
$ws.Range("E38").Value = '  -1.31%  '
$ws.Range("E39").Value = '  -4.89%  '
$ws.Range("D40").Value = '1.132.73'
$ws.Range("E40").Value = '  -0.70%  '
$ws.Range("D41").Value = '''6.203'
$ws.Range("E41").Value = '  -4.51%  '
$ws.Range("D42").Value = '''0.8755'
$ws.Range("E42").Value = '  -1.46%  '
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("D44").Value = '''99.52'
$ws.Range("E44").Value = '  -1.30%  '
$ws.Range("D45").Value = '1.891.73'
$ws.Range("E45").Value = '  -4.85%  '
$ws.Range("D46").Value = '''59.43'
$ws.Range("E46").Value = '  -3.88%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '''1.584'
$ws.Range("E47").Value = '  -1.71%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '''0.00000000107'
$ws.Range("E48").Value = '  -4.81%  '
$ws.Range("D49").Value = '''8.283'
$ws.Range("E49").Value = '  -2.44%  '
$ws.Range("D50").Value = '''0.05380'
$ws.Range("D51").Value = '''6.277'
$ws.Range("E51").Value = '  -1.75%  '
